$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ts ")

# --- Capture existing B-column contents (rows 6-10) before they get overwritten ---
$oldB6 = $ws.Range("B6").Value()
$oldB7 = $ws.Range("B7").Value()
$oldB8 = $ws.Range("B8").Value()
$oldB9 = $ws.Range("B9").Value()
$oldB10 = $ws.Range("B10").Value()

# --- Fill in the new "TEST CASE NO" counts in column D (rows 2-6) ---
$ws.Range("D2").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 1

# --- Row 11: newly populated with a fresh scenario (was previously blank) ---
$ws.Range("A11").Value = "filpkey_10"
$ws.Range("A11").Style = "Normal"
$ws.Range("A11").Font.Color = $ws.Range("A10").Font.Color
$ws.Range("B11").Value = $oldB10
$ws.Range("C11").Value = "positive"

# --- Shift the description/positive pairs for rows 6-10 down by one row ---
$ws.Range("B10").Value = $oldB9
$ws.Range("B9").Value = $oldB8
$ws.Range("B8").Value = $oldB7
$ws.Range("B7").Value = $oldB6
$ws.Range("C7").Value = "positive"

# --- Row 6 now holds a brand new scenario description ---
$ws.Range("B6").Value = "Validate moneychange symbol  funtionality"
$ws.Range("C6").Value = "positive"

# --- Formatting: the "wrap text" styling moves from row 6 to row 7 ---
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Style = "Normal"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Style = "Normal"
$ws.Range("B7").WrapText = $true
$ws.Range("B7").VerticalAlignment = -4160
$ws.Range("C7").VerticalAlignment = -4108

# --- Row heights: the taller row moves from row 6 to row 7 ---
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(7).RowHeight = 29

# --- Update the selection to match the saved view ---
$ws.Range("B7").Select()
